# Adds two new worksheets ("2303281656" and "2303281657") to the end of the
# workbook, each containing the department/institute/faculty abbreviation
# table with an "Obsolete terms" column. The first new sheet is the
# intermediate state where the Amsterdam Institute for Life and Environment
# dept. abbreviation was (incorrectly) set to "beta-aminlien"; the second new
# sheet reverts that back to "beta-aile" (already used in Yoda) and moves
# "beta-aminlien" into the Obsolete terms list instead.

$sheet9Data = @(
    @("Dept. abbr.", "Department", "Inst. abbr.", "Institute", "Fac. abbr.", "Faculty", "Obsolete terms"),
    @("fgb-acesh", "Academic Center for Education, Sport and Health", "vu-alab", "A-LAB", "acta", "Academic Centre for Dentistry Amsterdam", "acta-fda"),
    @("sbe-acc", "Accounting", "vu-aimms", "AIMMS", "fgb", "Faculty of Behavioural and Movement Sciences", "beta-ecsc"),
    @("beta-aminlien", "Amsterdam Institute for Life and Environment", "vu-abri", "Amsterdam Business Research Institute", "fgw", "Faculty of Humanities", "beta-mcb"),
    @("fgw-acha", "Art and Culture, History, Antiquity", "vu-ams", "Amsterdam Movement Sciences", "rch", "Faculty of Law", "vu-whocc"),
    @("beta-ai", "Athena Institute", "vu-an", "Amsterdam Neuroscience", "frt", "Faculty of Religion and Theology", "beta-aile"),
    @("frt-bp", "Beliefs and Practices", "vu-aph", "Amsterdam Public Health", "beta", "Faculty of Science", "beta-alife"),
    @("fgb-bp", "Biological Psychology", "vu-ard", "Amsterdam Reproduction & Development", "fsw", "Faculty of Social Sciences", "vu-cic"),
    @("beta-cncr", "Center for Neurogenomics and Cognitive Research", "vu-asi", "Amsterdam Sustainability Institute", "sbe", "School of Business and Economics", "vu-kcdi"),
    @("beta-cps", "Chemistry and Pharmaceutical Sciences", "vu-clue", "CLUE+", "gnk", "VUmc - School of Medical Sciences", "beta-eh"),
    @("fgb-cndp", "Clinical, Neuro- & Developmental Psychology", "vu-ibba", "IBBA", "", "", "vu-acwfs"),
    @("fsw-cosc", "Communication Science", "vu-isr", "Institute for Societal Resilience", "", "", "acta-acdia"),
    @("beta-cs", "Computer Science", "vu-ki", "Kooijmans Institute", "", "", ""),
    @("rch-cal", "Constitutional and Administrative Law", "vu-learn", "LEARN!", "", "", ""),
    @("rch-clc", "Criminal Law and Criminology", "vu-laser", "LaserLaB", "", "", ""),
    @("rch-dpl", "Dutch Private Law", "vu-ni", "Network Institute", "", "", ""),
    @("beta-es", "Earth Sciences", "vu-tain", "Talma Institute", "", "", ""),
    @("sbe-eds", "Econometrics and Data Science", "vu-ti", "Tinbergen Institute", "", "", ""),
    @("sbe-econ", "Economics", "", "", "", "", ""),
    @("fgb-efs", "Educational and Family Studies", "", "", "", "", ""),
    @("sbe-egs", "Ethics, Governance and Society", "", "", "", "", ""),
    @("fgb-eap", "Experimental and Applied Psychology", "", "", "", "", ""),
    @("sbe-fin", "Finance", "", "", "", "", ""),
    @("beta-hs", "Health Sciences", "", "", "", "", ""),
    @("beta-hsas", "History and Social Aspects of Science", "", "", "", "", ""),
    @("fgb-hms", "Human Movement Sciences", "", "", "", "", ""),
    @("beta-ies", "Institute for Environmental Studies", "", "", "", "", ""),
    @("beta-kggb", "Kars Group (Geo- and Bioarchaeology)", "", "", "", "", ""),
    @("sbe-kii", "Knowledge, Information and Innovation", "", "", "", "", ""),
    @("fgw-llc", "Language, Literature and Communication", "", "", "", "", ""),
    @("rch-ltlh", "Legal Theory and Legal History", "", "", "", "", ""),
    @("sbe-mo", "Management and Organisation", "", "", "", "", ""),
    @("sbe-mrk", "Marketing", "", "", "", "", ""),
    @("beta-math", "Mathematics", "", "", "", "", ""),
    @("rch-ntl", "Notary and Tax Law", "", "", "", "", ""),
    @("acta-owi", "OWI (ACTA)", "", "", "", "", ""),
    @("sbe-oa", "Operations Analytics", "", "", "", "", ""),
    @("acta-oii", "Oral Infections and Inflammation (OII)", "", "", "", "", ""),
    @("acta-orm", "Oral Regenerative Medicine (ORM)", "", "", "", "", ""),
    @("fsw-os", "Organization Sciences", "", "", "", "", ""),
    @("fgw-phil", "Philosophy", "", "", "", "", ""),
    @("beta-pa", "Physics and Astronomy", "", "", "", "", ""),
    @("fsw-pspa", "Political Science and Public Administration", "", "", "", "", ""),
    @("fsw-sca", "Social and Cultural Anthropology", "", "", "", "", ""),
    @("fsw-socio", "Sociology", "", "", "", "", ""),
    @("sbe-se", "Spatial Economics", "", "", "", "", ""),
    @("frt-tt", "Texts and Traditions", "", "", "", "", ""),
    @("rch-tls", "Transnational Legal Studies", "", "", "", "", ""),
    @("sbe-vsee", "VU SBE Executive Education", "", "", "", "", "")
)

$sheet10Data = @(
    @("Dept. abbr.", "Department", "Inst. abbr.", "Institute", "Fac. abbr.", "Faculty", "Obsolete terms"),
    @("fgb-acesh", "Academic Center for Education, Sport and Health", "vu-alab", "A-LAB", "acta", "Academic Centre for Dentistry Amsterdam", "acta-fda"),
    @("sbe-acc", "Accounting", "vu-aimms", "AIMMS", "fgb", "Faculty of Behavioural and Movement Sciences", "acta-acdia"),
    @("beta-aile", "Amsterdam Institute for Life and Environment", "vu-abri", "Amsterdam Business Research Institute", "fgw", "Faculty of Humanities", "vu-whocc"),
    @("fgw-acha", "Art and Culture, History, Antiquity", "vu-ams", "Amsterdam Movement Sciences", "rch", "Faculty of Law", "beta-eh"),
    @("beta-ai", "Athena Institute", "vu-an", "Amsterdam Neuroscience", "frt", "Faculty of Religion and Theology", "beta-aminlien"),
    @("frt-bp", "Beliefs and Practices", "vu-aph", "Amsterdam Public Health", "beta", "Faculty of Science", "beta-mcb"),
    @("fgb-bp", "Biological Psychology", "vu-ard", "Amsterdam Reproduction & Development", "fsw", "Faculty of Social Sciences", "beta-ecsc"),
    @("beta-cncr", "Center for Neurogenomics and Cognitive Research", "vu-asi", "Amsterdam Sustainability Institute", "sbe", "School of Business and Economics", "vu-acwfs"),
    @("beta-cps", "Chemistry and Pharmaceutical Sciences", "vu-clue", "CLUE+", "gnk", "VUmc - School of Medical Sciences", "vu-cic"),
    @("fgb-cndp", "Clinical, Neuro- & Developmental Psychology", "vu-ibba", "IBBA", "", "", "vu-kcdi"),
    @("fsw-cosc", "Communication Science", "vu-isr", "Institute for Societal Resilience", "", "", ""),
    @("beta-cs", "Computer Science", "vu-ki", "Kooijmans Institute", "", "", ""),
    @("rch-cal", "Constitutional and Administrative Law", "vu-learn", "LEARN!", "", "", ""),
    @("rch-clc", "Criminal Law and Criminology", "vu-laser", "LaserLaB", "", "", ""),
    @("rch-dpl", "Dutch Private Law", "vu-ni", "Network Institute", "", "", ""),
    @("beta-es", "Earth Sciences", "vu-tain", "Talma Institute", "", "", ""),
    @("sbe-eds", "Econometrics and Data Science", "vu-ti", "Tinbergen Institute", "", "", ""),
    @("sbe-econ", "Economics", "", "", "", "", ""),
    @("fgb-efs", "Educational and Family Studies", "", "", "", "", ""),
    @("sbe-egs", "Ethics, Governance and Society", "", "", "", "", ""),
    @("fgb-eap", "Experimental and Applied Psychology", "", "", "", "", ""),
    @("sbe-fin", "Finance", "", "", "", "", ""),
    @("beta-hs", "Health Sciences", "", "", "", "", ""),
    @("beta-hsas", "History and Social Aspects of Science", "", "", "", "", ""),
    @("fgb-hms", "Human Movement Sciences", "", "", "", "", ""),
    @("beta-ies", "Institute for Environmental Studies", "", "", "", "", ""),
    @("beta-kggb", "Kars Group (Geo- and Bioarchaeology)", "", "", "", "", ""),
    @("sbe-kii", "Knowledge, Information and Innovation", "", "", "", "", ""),
    @("fgw-llc", "Language, Literature and Communication", "", "", "", "", ""),
    @("rch-ltlh", "Legal Theory and Legal History", "", "", "", "", ""),
    @("sbe-mo", "Management and Organisation", "", "", "", "", ""),
    @("sbe-mrk", "Marketing", "", "", "", "", ""),
    @("beta-math", "Mathematics", "", "", "", "", ""),
    @("rch-ntl", "Notary and Tax Law", "", "", "", "", ""),
    @("acta-owi", "OWI (ACTA)", "", "", "", "", ""),
    @("sbe-oa", "Operations Analytics", "", "", "", "", ""),
    @("acta-oii", "Oral Infections and Inflammation (OII)", "", "", "", "", ""),
    @("acta-orm", "Oral Regenerative Medicine (ORM)", "", "", "", "", ""),
    @("fsw-os", "Organization Sciences", "", "", "", "", ""),
    @("fgw-phil", "Philosophy", "", "", "", "", ""),
    @("beta-pa", "Physics and Astronomy", "", "", "", "", ""),
    @("fsw-pspa", "Political Science and Public Administration", "", "", "", "", ""),
    @("fsw-sca", "Social and Cultural Anthropology", "", "", "", "", ""),
    @("fsw-socio", "Sociology", "", "", "", "", ""),
    @("sbe-se", "Spatial Economics", "", "", "", "", ""),
    @("frt-tt", "Texts and Traditions", "", "", "", "", ""),
    @("rch-tls", "Transnational Legal Studies", "", "", "", "", ""),
    @("sbe-vsee", "VU SBE Executive Education", "", "", "", "", "")
)

$wb = $excel.ActiveWorkbook

# Excel stores column widths padded by ~5/6 of a character (the default-font
# glyph-padding quirk), so to land on an exact target "width" in the saved
# xlsx we have to pre-subtract that offset before assigning ColumnWidth.
$colWidthOffset = 5.0 / 6.0

function Add-AbbreviationSheet {
    param([string]$SheetName, [array]$RowsData, [array]$ColWidths)

    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newSheet = $wb.Worksheets.Add($null, $lastSheet)
    $newSheet.Name = $SheetName

    for ($i = 0; $i -lt $ColWidths.Length; $i++) {
        $newSheet.Columns.Item($i + 1).ColumnWidth = $ColWidths[$i] - $colWidthOffset
    }

    $r = 1
    foreach ($row in $RowsData) {
        $c = 1
        foreach ($val in $row) {
            if ($val -ne "") {
                $cell = $newSheet.Cells.Item($r, $c)
                $cell.Value = $val
                if ($r -eq 1) {
                    $cell.Font.Bold = $true
                }
            }
            $c = $c + 1
        }
        $r = $r + 1
    }

    $newSheet.Range("A1").Select()
}

Add-AbbreviationSheet "2303281656" $sheet9Data @(13, 47, 11, 37, 10, 44, 14)
Add-AbbreviationSheet "2303281657" $sheet10Data @(11, 47, 11, 37, 10, 44, 14)

Write-Host "Added sheets, total count:" $wb.Worksheets.Count
